$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Authentication / login
$ws.Range("A5").Value = "Authentication"
$ws.Range("B5").Value = "post"
$ws.Range("C5").Value = "/authentication"
$ws.Range("F5").Value = "email, password"
$ws.Range("G5").Value = "login"
$ws.Range("H5").Value = "token, user info"

# Row 6 - signup
$ws.Range("B6").Value = "post"
$ws.Range("C6").Value = "/signup"

# F6: rich text body for signup - required fields (bold "(not null)")
$f6text = "email, password, name, type (Instructor, Administrator, Student) (not null)"
$ws.Range("F6").Value = $f6text
$f6BoldStart = $f6text.IndexOf("(not null)") + 1
$f6BoldLen = "(not null)".Length
$ws.Range("F6").Characters($f6BoldStart, $f6BoldLen).Font.Bold = $true

# Row 7 - optional personal fields (bold "(option)")
$f7text = "first_name, last_name, gender, birthday, avatar_url, user_address (option)"
$ws.Range("F7").Value = $f7text
$f7BoldStart = $f7text.IndexOf("(option)") + 1
$f7BoldLen = "(option)".Length
$ws.Range("F7").Characters($f7BoldStart, $f7BoldLen).Font.Bold = $true

# Row 8 - optional instructor fields (bold "(option)")
$f8text = "if type is Instructor: job_title, short_description, full_description (option)"
$ws.Range("F8").Value = $f8text
$f8BoldStart = $f8text.IndexOf("(option)") + 1
$f8BoldLen = "(option)".Length
$ws.Range("F8").Characters($f8BoldStart, $f8BoldLen).Font.Bold = $true

# G6 - signup (filled last, after F7/F8, to match shared string order)
$ws.Range("G6").Value = "signup"

# Selection as in target (activeCell G9 sqref G9)
$ws.Range("G9").Select() | Out-Null

# Column width adjustments to reflect new (longer) content - mirrors the
# bestFit auto-sizing Excel performs after the new rows/text are added.
$ws.Columns.Item(1).ColumnWidth = 13.592447916666666
$ws.Columns.Item(4).ColumnWidth = 5.307291666666667
$ws.Columns.Item(6).ColumnWidth = 68.02213541666667
$ws.Columns.Item(8).ColumnWidth = 14.307291666666666
